$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Dashboard")

# Test case 4 is now partly finished: the first two automated test cases
# (rows 2 & 3) passed, and the third (row 4) failed. All four were
# previously left as "Skipped".
$ws.Range("D2").Value = "Pass"
$ws.Range("D3").Value = "Pass"
$ws.Range("D4").Value = "Failed"

# Leave the selection where the author left off.
$ws.Range("B2").Select() | Out-Null
